# Add a new "latest_reported_idrc" column into the "Summary" sheet,
# shifting the existing cost22/cost23/cost24 columns one to the right
# and renaming them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Insert a new blank column at C, pushing old C/D/E -> D/E/F
# (this also carries over the header cell style from the old column C)
$ws.Columns.Item(3).Insert()

# --- Header row ---
$ws.Range("C1").Value = "latest_reported_idrc"
$ws.Range("D1").Value = "additional_cost_2022"
$ws.Range("E1").Value = "additional_cost_2023"
$ws.Range("F1").Value = "additional_cost_2024"

# --- New "latest_reported_idrc" values (column C), by row ---
$idrcValues = @{
    2  = 63.12
    3  = 250.67
    4  = 463.6
    5  = 368.35
    6  = 6.42
    7  = 2732.06
    8  = 63.2
    9  = 232.08
    10 = 68.77
    11 = 1156.48
    12 = 1447.37
    13 = 36.55
    14 = 1.38
    15 = 50.29
    16 = 4.33
    17 = 556.1
    18 = 0.27
    20 = 407.65
    21 = 52.25
    22 = 11.67
    23 = 16.88
    24 = 10.08
    25 = 1.21
    26 = 2.33
    27 = 87.91
    28 = 4745.18
}

foreach ($row in $idrcValues.Keys) {
    $ws.Cells.Item($row, 3).Value = $idrcValues[$row]
}

# Row 19 (Luxembourg) has no reported value for latest_reported_idrc,
# leave the cell blank.
$ws.Cells.Item(19, 3).Value = ""

Write-Output "done"
